$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (K1:P1) - swap label order to +/- pairs
# (order of assignment matters for shared-string table ordering)
$ws.Range("K1").Value = "``+ Nerves"
$ws.Range("L1").Value = "``- Nerves"
$ws.Range("M1").Value = "``+ Aggro"
$ws.Range("O1").Value = "``+ Fatigue"
$ws.Range("N1").Value = "``- Aggro"
$ws.Range("P1").Value = "``- Fatigue"

# Row labels in column J
$ws.Range("J2").Value = "``+ Nerves"
$ws.Range("J3").Value = "``- Nerves"
$ws.Range("J4").Value = "``+ Aggro"
$ws.Range("J5").Value = "``- Aggro"
$ws.Range("J6").Value = "``+ Fatigue"
$ws.Range("J7").Value = "``- Fatigue"

# New data cells
$ws.Range("O3").Value = "Speculate"
$ws.Range("O4").Value = "Coax"

# Update selection
$ws.Range("P3").Select()
